$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 284, shifting existing rows 284-339 down to 285-340
$ws.Rows.Item(284).Insert()

# Populate the newly inserted row 284 with the new data point
$ws.Cells.Item(284, 1).Value = 3
$ws.Cells.Item(284, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(284, 3).Value = "Coquimbo"
$ws.Cells.Item(284, 4).Value = 45258
$ws.Cells.Item(284, 5).Value = 5
$ws.Cells.Item(284, 6).Value = 100112026
$ws.Cells.Item(284, 7).Value = "Haba"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 70
$ws.Cells.Item(284, 11).Value = 10000
$ws.Cells.Item(284, 12).Value = 11000
$ws.Cells.Item(284, 13).Value = 10500
$ws.Cells.Item(284, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(284, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(284, 16).Value = 420
$ws.Cells.Item(284, 17).Value = 25
$ws.Cells.Item(284, 18).Value = "Hortaliza"
